$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.427.39"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "1.850.30"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.37"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6297"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07650"
$ws.Range("E8").Value = "  +0.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2914"
$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.92"
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("D11").Value = "2.131.99"
$ws.Range("E11").Value = "  +15.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07746"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.034"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6813"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001066"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.38"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "29.501.63"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.58"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.34"
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.468"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.95"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.447"
$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.70"
$ws.Range("E27").Value = "  +0.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.378"
$ws.Range("E28").Value = "  +6.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.464"
$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05624"
$ws.Range("E30").Value = "  +0.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.133"
$ws.Range("E31").Value = "  +0.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.054"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.846"
$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.165"
$ws.Range("E34").Value = "  +0.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7025"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.593"
$ws.Range("E36").Value = "  +0.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01805"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").Value = "1.229.39"
$ws.Range("E38").Value = "  -0.85%  "

$ws.Range("E39").Value = "  -2.10%  "

$ws.Range("E40").Value = "  +0.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9093"
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.45"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("E45").Value = "  -1.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.196"

$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("E48").Value = "  +3.18%  "

$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.983"
$ws.Range("E50").Value = "  -0.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05710"
$ws.Range("E51").Value = "  +0.11%  "
